$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 10 (weekly update: new week's prices pushed in,
# existing rows shift down by 4). Excel will inherit the formatting (incl. the
# date style on column D) from the row above, exactly like interactive usage.
$ws.Range("A10:A13").EntireRow.Insert()

# Row 10: Papaya, Especial
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 45118
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100108
$ws.Cells.Item(10, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(10, 9).Value = 100108004
$ws.Cells.Item(10, 10).Value = "Papaya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 140
$ws.Cells.Item(10, 14).Value = 24000
$ws.Cells.Item(10, 15).Value = 24000
$ws.Cells.Item(10, 16).Value = 24000
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 19).Value = 2400
$ws.Cells.Item(10, 20).Value = 10

# Row 11: Papaya, Primera
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 45118
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100108
$ws.Cells.Item(11, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(11, 9).Value = 100108004
$ws.Cells.Item(11, 10).Value = "Papaya"
$ws.Cells.Item(11, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 160
$ws.Cells.Item(11, 14).Value = 20000
$ws.Cells.Item(11, 15).Value = 20000
$ws.Cells.Item(11, 16).Value = 20000
$ws.Cells.Item(11, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 19).Value = 2000
$ws.Cells.Item(11, 20).Value = 10

# Row 12: Papaya, Segunda
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 45118
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100108
$ws.Cells.Item(12, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(12, 9).Value = 100108004
$ws.Cells.Item(12, 10).Value = "Papaya"
$ws.Cells.Item(12, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 180
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 15000
$ws.Cells.Item(12, 16).Value = 15000
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 19).Value = 1500
$ws.Cells.Item(12, 20).Value = 10

# Row 13: Papaya, Tercera
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 45118
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100108
$ws.Cells.Item(13, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(13, 9).Value = 100108004
$ws.Cells.Item(13, 10).Value = "Papaya"
$ws.Cells.Item(13, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(13, 12).Value = "Tercera"
$ws.Cells.Item(13, 13).Value = 75
$ws.Cells.Item(13, 14).Value = 12000
$ws.Cells.Item(13, 15).Value = 12000
$ws.Cells.Item(13, 16).Value = 12000
$ws.Cells.Item(13, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(13, 19).Value = 1200
$ws.Cells.Item(13, 20).Value = 10
